$wb = $excel.ActiveWorkbook

$newTime = "03:55:07"

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 13"

$ws1.Range("A16").Value = $newTime
$ws1.Range("B16").Value = "05:35"
$ws1.Range("C16").Value = "215B_EL PATO"
$ws1.Range("D16").Value = 100
$ws1.Range("E16").Value = "LP1912"

$ws1.Range("A17").Value = $newTime
$ws1.Range("B17").Value = "05:46"
$ws1.Range("C17").Value = "15_ABASTO"
$ws1.Range("D17").Value = 111
$ws1.Range("E17").Value = "LP1912"

$ws1.Range("A18").Value = $newTime
$ws1.Range("B18").Value = "05:54"
$ws1.Range("C18").Value = "10_OLMOS"
$ws1.Range("D18").Value = 119
$ws1.Range("E18").Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 5"

$ws2.Range("A10").Value = $newTime
$ws2.Range("B10").Value = "05:35"
$ws2.Range("C10").Value = "215B_EL PATO"
$ws2.Range("D10").Value = 100
$ws2.Range("E10").Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newTime"
$ws3.Range("A3").Value = "Total filas: 2"

$ws3.Range("A7").Value = $newTime
$ws3.Range("B7").Value = "05:44"
$ws3.Range("C7").Value = "215A_LA PLATA"
$ws3.Range("D7").Value = 109
$ws3.Range("E7").Value = "L6173"
